# Update the "Förändrad" (Changed) date column (C) for all data rows.
# The workbook tracks a refresh/changed timestamp in column C for every
# record; this automatic update bumps that date serial from 45186 to 45188
# (2023-09-17 -> 2023-09-19) for every data row on the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$lastRow = $ws.Cells.Item($ws.Rows.Count, 1).End(-4162).Row
if ($lastRow -lt 2) {
    $lastRow = $ws.UsedRange.Rows.Count
}

for ($r = 2; $r -le $lastRow; $r++) {
    $cell = $ws.Cells.Item($r, 3)
    if ($cell.Value2 -eq 45186) {
        $cell.Value2 = 45188
    }
}
